# Append the "Results of live testing analysis" section after the last
# table, right before the document's closing section properties.
#
# The document currently ends with:
#   ... </w:tbl> <w:p/> <w:sectPr> ... </w:sectPr>
#
# and the target state is:
#   ... </w:tbl> <w:p/> <new paragraphs...> <w:p/> <w:sectPr> ... </w:sectPr>
#
# The existing trailing empty paragraph (right after the table) must be
# left completely untouched, so new content is added via a fresh paragraph
# inserted after it rather than by rewriting that paragraph in place.

$d = $word.ActiveDocument

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$newXml = @"
<w:p $wns>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>Results of live testing analysis</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t xml:space="preserve"> – Q1 + Q2</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Try to explain what you have found.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>What was the analysis from the strategy used?</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> PnL, Win Rate, Avg Win/Loss etc… average price over time (weekly) and how this effects </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>pnl</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or win rate.</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Was there anything you found that was unexpected that may improve the system?</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:r>
    <w:t xml:space="preserve"> I.E. now a win rate has been established, a more accurate risk/reward ratio can be applied. Using this strategy, or indeed any sports betting strategy, all that really matters is the odds or price you enter the market. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">It is yet unconfirmed if the stats used for the strategy have any correlation with the price used. </w:t>
  </w:r>
  <w:r>
    <w:t>As football results can appear random the odds reflect the results, no match will come in just because you want it to. Which seems self-explanatory but when you want something to happen it feels wrong when it doesn’t. Therefore</w:t>
  </w:r>
  <w:r>
    <w:t>,</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> I theorise that a lower average price must be achieved </w:t>
  </w:r>
  <w:r>
    <w:t>to</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>produce</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> a long term positive expected value (+EV). I think that from the results of Q1 + 2, albeit fragmented, it shows that a lower price entered on the lay side over time results in +EV. This could be proof </w:t>
  </w:r>
  <w:r>
    <w:t>of edge this season without the need to risk unbalanced amounts when the bet is in my favour.</w:t>
  </w:r>
</w:p>
<w:p $wns/>
"@

$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)

# Insert a fresh paragraph after the existing (final, empty) paragraph so the
# latter's identity/attributes are left completely untouched.
$lastPara.Range.InsertParagraphAfter()

# The freshly-inserted paragraph is now the new last paragraph; collapse its
# range to the end and swap in the full run of new paragraphs via InsertXML
# (inserting XML into a collapsed, empty trailing paragraph replaces that
# paragraph with the supplied content instead of leaving a stray empty one).
$n2 = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($n2).Range
$target.Collapse(0)
$target.InsertXML($newXml)
